$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Cells.Item(32, 8).Value = 1021.3125
$ws.Cells.Item(32, 9).Value = 883.1667
$ws.Cells.Item(32, 10).Value = 1104.2
$ws.Cells.Item(32, 11).Value = 883.1667
$ws.Cells.Item(32, 12).Value = 1104.2
$ws.Cells.Item(32, 13).Value = -557.1667
$ws.Cells.Item(32, 14).Value = -1756.2
# Row 62
$ws.Cells.Item(62, 8).Value = 1639.5
$ws.Cells.Item(62, 9).Value = 2139.0833
$ws.Cells.Item(62, 10).Value = 1040
$ws.Cells.Item(62, 11).Value = 2139.0833
$ws.Cells.Item(62, 12).Value = 1040
$ws.Cells.Item(62, 13).Value = -1515.0833
$ws.Cells.Item(62, 14).Value = -2288
# Row 64
$ws.Cells.Item(64, 8).Value = 3105.7144
$ws.Cells.Item(64, 9).Value = 3128.3572
$ws.Cells.Item(64, 10).Value = 3060.4285
$ws.Cells.Item(64, 11).Value = 3128.3572
$ws.Cells.Item(64, 12).Value = 3060.4285
$ws.Cells.Item(64, 13).Value = -2880.3572
$ws.Cells.Item(64, 14).Value = -3556.4285
# Row 65
$ws.Cells.Item(65, 8).Value = 1639.5
$ws.Cells.Item(65, 9).Value = 2139.0833
$ws.Cells.Item(65, 10).Value = 1040
$ws.Cells.Item(65, 11).Value = 10695.4165
$ws.Cells.Item(65, 12).Value = 5200
$ws.Cells.Item(65, 13).Value = -7575.416499999999
$ws.Cells.Item(65, 14).Value = -11440
# Row 67
$ws.Cells.Item(67, 8).Value = 3105.7144
$ws.Cells.Item(67, 9).Value = 3128.3572
$ws.Cells.Item(67, 10).Value = 3060.4285
$ws.Cells.Item(67, 11).Value = 3128.3572
$ws.Cells.Item(67, 12).Value = 3060.4285
$ws.Cells.Item(67, 13).Value = -2270.3572
$ws.Cells.Item(67, 14).Value = -4776.4285
# Row 98
$ws.Cells.Item(98, 8).Value = 1839.6875
$ws.Cells.Item(98, 9).Value = 1787.9166
$ws.Cells.Item(98, 10).Value = 1995
$ws.Cells.Item(98, 11).Value = 1787.9166
$ws.Cells.Item(98, 12).Value = 1995
$ws.Cells.Item(98, 13).Value = -289.9166
$ws.Cells.Item(98, 14).Value = -4991
# Row 122
$ws.Cells.Item(122, 8).Value = 1839.6875
$ws.Cells.Item(122, 9).Value = 1787.9166
$ws.Cells.Item(122, 10).Value = 1995
$ws.Cells.Item(122, 11).Value = 5363.7498
$ws.Cells.Item(122, 12).Value = 5985
$ws.Cells.Item(122, 13).Value = -2913.7498
$ws.Cells.Item(122, 14).Value = -10885
# Row 132
$ws.Cells.Item(132, 8).Value = 4723063
$ws.Cells.Item(132, 9).Value = 4923.026
$ws.Cells.Item(132, 10).Value = 17866454
$ws.Cells.Item(132, 11).Value = 14769.078
$ws.Cells.Item(132, 12).Value = 53599362
$ws.Cells.Item(132, 13).Value = -12239.078
$ws.Cells.Item(132, 14).Value = -53604422

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 14691.789
$ws.Cells.Item(32, 9).Value = 19566.291
$ws.Cells.Item(32, 10).Value = 6335.5
$ws.Cells.Item(32, 11).Value = 19566.291
$ws.Cells.Item(32, 12).Value = 6335.5
$ws.Cells.Item(32, 13).Value = -19279.291
$ws.Cells.Item(32, 14).Value = -6909.5
# Row 61
$ws.Cells.Item(61, 8).Value = 9617250
$ws.Cells.Item(61, 9).Value = 13159787
$ws.Cells.Item(61, 10).Value = 1793.8572
$ws.Cells.Item(61, 11).Value = 13159787
$ws.Cells.Item(61, 12).Value = 1793.8572
$ws.Cells.Item(61, 13).Value = -13159575
$ws.Cells.Item(61, 14).Value = -2217.8572
# Row 74
$ws.Cells.Item(74, 8).Value = 7937888
$ws.Cells.Item(74, 9).Value = 9434929
$ws.Cells.Item(74, 10).Value = 3569
$ws.Cells.Item(74, 11).Value = 9434929
$ws.Cells.Item(74, 12).Value = 3569
$ws.Cells.Item(74, 13).Value = -9434055
$ws.Cells.Item(74, 14).Value = -5317
# Row 77
$ws.Cells.Item(77, 8).Value = 7937888
$ws.Cells.Item(77, 9).Value = 9434929
$ws.Cells.Item(77, 10).Value = 3569
$ws.Cells.Item(77, 11).Value = 47174645
$ws.Cells.Item(77, 12).Value = 17845
$ws.Cells.Item(77, 13).Value = -47170277
$ws.Cells.Item(77, 14).Value = -26581
# Row 102
$ws.Cells.Item(102, 8).Value = 11285.714
$ws.Cells.Item(102, 9).Value = 10666.667
$ws.Cells.Item(102, 10).Value = 15000
$ws.Cells.Item(102, 11).Value = 10666.667
$ws.Cells.Item(102, 12).Value = 15000
$ws.Cells.Item(102, 13).Value = -9044.666999999999
$ws.Cells.Item(102, 14).Value = -18244
# Row 136
$ws.Cells.Item(136, 8).Value = 9617250
$ws.Cells.Item(136, 9).Value = 13159787
$ws.Cells.Item(136, 10).Value = 1793.8572
$ws.Cells.Item(136, 11).Value = 39479361
$ws.Cells.Item(136, 12).Value = 5381.571599999999
$ws.Cells.Item(136, 13).Value = -39476811
$ws.Cells.Item(136, 14).Value = -10481.5716

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 5407.224
$ws.Cells.Item(134, 9).Value = 4642.9287
$ws.Cells.Item(134, 10).Value = 7413.5
$ws.Cells.Item(134, 11).Value = 13928.7861
$ws.Cells.Item(134, 12).Value = 22240.5
$ws.Cells.Item(134, 13).Value = -11393.7861
$ws.Cells.Item(134, 14).Value = -27310.5

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 6176710
$ws.Cells.Item(31, 9).Value = 5734.067
$ws.Cells.Item(31, 10).Value = 13890430
$ws.Cells.Item(31, 11).Value = 5734.067
$ws.Cells.Item(31, 12).Value = 13890430
$ws.Cells.Item(31, 13).Value = -5439.067
$ws.Cells.Item(31, 14).Value = -13891020
# Row 34
$ws.Cells.Item(34, 8).Value = 6176710
$ws.Cells.Item(34, 9).Value = 5734.067
$ws.Cells.Item(34, 10).Value = 13890430
$ws.Cells.Item(34, 11).Value = 5734.067
$ws.Cells.Item(34, 12).Value = 13890430
$ws.Cells.Item(34, 13).Value = -5532.067
$ws.Cells.Item(34, 14).Value = -13890834
# Row 58
$ws.Cells.Item(58, 8).Value = 2307.7878
$ws.Cells.Item(58, 9).Value = 1135.125
$ws.Cells.Item(58, 10).Value = 3411.4707
$ws.Cells.Item(58, 11).Value = 1135.125
$ws.Cells.Item(58, 12).Value = 3411.4707
$ws.Cells.Item(58, 13).Value = -932.125
$ws.Cells.Item(58, 14).Value = -3817.4707
# Row 132
$ws.Cells.Item(132, 8).Value = 2139.4102
$ws.Cells.Item(132, 9).Value = 1709.5588
$ws.Cells.Item(132, 10).Value = 5062.4
$ws.Cells.Item(132, 11).Value = 5128.6764
$ws.Cells.Item(132, 12).Value = 15187.2
$ws.Cells.Item(132, 13).Value = -2598.6764
$ws.Cells.Item(132, 14).Value = -20247.2
# Row 136
$ws.Cells.Item(136, 8).Value = 2307.7878
$ws.Cells.Item(136, 9).Value = 1135.125
$ws.Cells.Item(136, 10).Value = 3411.4707
$ws.Cells.Item(136, 11).Value = 3405.375
$ws.Cells.Item(136, 12).Value = 10234.4121
$ws.Cells.Item(136, 13).Value = -855.375
$ws.Cells.Item(136, 14).Value = -15334.4121

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 74
$ws.Cells.Item(74, 8).Value = 7124.875
$ws.Cells.Item(74, 9).Value = 3140
$ws.Cells.Item(74, 10).Value = 13766.333
$ws.Cells.Item(74, 11).Value = 9420
$ws.Cells.Item(74, 12).Value = 41298.999
$ws.Cells.Item(74, 13).Value = -8359
$ws.Cells.Item(74, 14).Value = -43420.999
# Row 77
$ws.Cells.Item(77, 8).Value = 7124.875
$ws.Cells.Item(77, 9).Value = 3140
$ws.Cells.Item(77, 10).Value = 13766.333
$ws.Cells.Item(77, 11).Value = 28260
$ws.Cells.Item(77, 12).Value = 123896.997
$ws.Cells.Item(77, 13).Value = -22956
$ws.Cells.Item(77, 14).Value = -134504.997
# Row 113
$ws.Cells.Item(113, 8).Value = 1168.2565
$ws.Cells.Item(113, 9).Value = 1123.2858
$ws.Cells.Item(113, 10).Value = 1193.44
$ws.Cells.Item(113, 11).Value = 3369.8574
$ws.Cells.Item(113, 12).Value = 3580.32
$ws.Cells.Item(113, 13).Value = -1199.8574
$ws.Cells.Item(113, 14).Value = -7920.32
# Row 131
$ws.Cells.Item(131, 8).Value = 845.52
$ws.Cells.Item(131, 9).Value = 266.66666
$ws.Cells.Item(131, 10).Value = 902.7692
$ws.Cells.Item(131, 11).Value = 799.9999799999999
$ws.Cells.Item(131, 12).Value = 2708.3076
$ws.Cells.Item(131, 13).Value = 4240.00002
$ws.Cells.Item(131, 14).Value = -12788.3076
# Row 132
$ws.Cells.Item(132, 8).Value = 2128.2188
$ws.Cells.Item(132, 9).Value = 919
$ws.Cells.Item(132, 10).Value = 3337.4375
$ws.Cells.Item(132, 11).Value = 8271
$ws.Cells.Item(132, 12).Value = 30036.9375
$ws.Cells.Item(132, 13).Value = -5741
$ws.Cells.Item(132, 14).Value = -35096.9375

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 102
$ws.Cells.Item(102, 8).Value = 45467.5
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 45467.5
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 12).Value = 45467.5
$ws.Cells.Item(102, 14).Value = -51957.5
# Row 132
$ws.Cells.Item(132, 8).Value = 10173.182
$ws.Cells.Item(132, 9).Value = 6933.4443
$ws.Cells.Item(132, 10).Value = 14060.866
$ws.Cells.Item(132, 11).Value = 20800.3329
$ws.Cells.Item(132, 12).Value = 42182.598
$ws.Cells.Item(132, 13).Value = -18270.3329
$ws.Cells.Item(132, 14).Value = -47242.598
# Row 136
$ws.Cells.Item(136, 8).Value = 34889200
$ws.Cells.Item(136, 9).Value = 46876988
$ws.Cells.Item(136, 10).Value = 15631.818
$ws.Cells.Item(136, 11).Value = 140630964
$ws.Cells.Item(136, 12).Value = 46895.454
$ws.Cells.Item(136, 13).Value = -140628414
$ws.Cells.Item(136, 14).Value = -51995.454

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 5548.067
$ws.Cells.Item(132, 9).Value = 6898.65
$ws.Cells.Item(132, 10).Value = 2846.9
$ws.Cells.Item(132, 11).Value = 20695.95
$ws.Cells.Item(132, 12).Value = 8540.700000000001
$ws.Cells.Item(132, 13).Value = -18165.95
$ws.Cells.Item(132, 14).Value = -13600.7
# Row 136
$ws.Cells.Item(136, 8).Value = 1182.7826
$ws.Cells.Item(136, 9).Value = 1100.409
$ws.Cells.Item(136, 10).Value = 2995
$ws.Cells.Item(136, 11).Value = 3301.227
$ws.Cells.Item(136, 12).Value = 8985
$ws.Cells.Item(136, 13).Value = -751.2270000000003
$ws.Cells.Item(136, 14).Value = -14085
